$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ngf"
$ws.Range("C2").Value = "Ntrk1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1930666666666666
$ws.Range("H2").Value = 0.5791999999999999
$ws.Range("I2").Value = 0.01292026122037801
$ws.Range("J2").Value = 0.01292026122037801
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1211523333333333
$ws.Range("N2").Value = 0.363457
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.02339047715555555
$ws.Range("R2").Value = 0.2105142944
$ws.Range("S2").Value = 0.01292026122037801
$ws.Range("T2").Value = 0.01292026122037801

# --- Row 3 ---
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ngf"
$ws.Range("C3").Value = "Ntrk1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.666149
$ws.Range("H3").Value = 1.998447
$ws.Range("I3").Value = 0.04457951877603724
$ws.Range("J3").Value = 0.04457951877603725
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.1211523333333333
$ws.Range("N3").Value = 0.363457
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.08070550569766667
$ws.Range("R3").Value = 0.7263495512790001
$ws.Range("S3").Value = 0.04457951877603724
$ws.Range("T3").Value = 0.04457951877603725

# --- Row 4 ---
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Ngf"
$ws.Range("C4").Value = "Ntrk1"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 14.08372266666667
$ws.Range("H4").Value = 42.251168
$ws.Range("I4").Value = 0.9425002200035847
$ws.Range("J4").Value = 0.9425002200035848
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1211523333333333
$ws.Range("N4").Value = 0.363457
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 1.706275863086222
$ws.Range("R4").Value = 15.356482767776
$ws.Range("S4").Value = 0.9425002200035847
$ws.Range("T4").Value = 0.9425002200035848

# --- Remove rows 5-7 (no longer needed) ---
$ws.Range("A5:T7").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
